$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)  # 2024-Match

$ws.Range('K8').Value = 'Eric Wiberg'

$ws.Range('J9').Value = 'Eric Wiberg'

$ws.Range('K16').Value = 'Eric Wiberg'

$ws.Range('J17').Value = 'Eric Wiberg'

$ws.Range('K20').Value = 'Eric Wiberg'

$ws.Range('J21').Value = 'Eric Wiberg'

$ws.Range('K24').Value = 'Eric Wiberg'

$ws.Range('J25').Value = 'Eric Wiberg'

$ws.Range('K28').Value = 'Eric Wiberg'

$ws.Range('J29').Value = 'Eric Wiberg'

$ws.Range('H30').Value = 2
$ws.Range('I30').Value = 1
$ws.Range('J30').Value = 'Dean Zoulamis'
$ws.Range('K30').Value = 'Patrik Udeh'

$ws.Range('H31').Value = 1
$ws.Range('I31').Value = 1
$ws.Range('J31').Value = 'Patrik Udeh'
$ws.Range('K31').Value = 'Dean Zoulamis'

$ws.Range('H32').Value = 2
$ws.Range('I32').Value = 3
$ws.Range('J32').Value = 'Caleb Lamb'
$ws.Range('K32').Value = 'Eric Wiberg'

$ws.Range('H33').Value = 1
$ws.Range('I33').Value = 6
$ws.Range('J33').Value = 'Eric Wiberg'
$ws.Range('K33').Value = 'Caleb Lamb'

$ws.Range('H34').Value = 3
$ws.Range('I34').Value = 1
$ws.Range('J34').Value = 'Dean Zoulamis'
$ws.Range('K34').Value = 'Patrik Udeh'

$ws.Range('H35').Value = 4
$ws.Range('I35').Value = 1
$ws.Range('J35').Value = 'Patrik Udeh'
$ws.Range('K35').Value = 'Dean Zoulamis'

$ws.Range('H36').Value = 1
$ws.Range('I36').Value = 2
$ws.Range('J36').Value = 'Patrik Udeh'
$ws.Range('K36').Value = 'Eric Wiberg'

$ws.Range('H37').Value = 0
$ws.Range('I37').Value = 6
$ws.Range('J37').Value = 'Eric Wiberg'
$ws.Range('K37').Value = 'Patrik Udeh'

$ws.Range('H38').Value = 7
$ws.Range('I38').Value = 1
$ws.Range('J38').Value = 'Dean Zoulamis'
$ws.Range('K38').Value = 'Patrik Udeh'

$ws.Range('H39').Value = 0
$ws.Range('I39').Value = 4
$ws.Range('J39').Value = 'Patrik Udeh'
$ws.Range('K39').Value = 'Dean Zoulamis'

$ws.Range('H40').Value = 0
$ws.Range('I40').Value = 4
$ws.Range('J40').Value = 'Eric Wiberg'
$ws.Range('K40').Value = 'Patrik Udeh'

$ws.Range('H41').Value = 0
$ws.Range('I41').Value = 0
$ws.Range('J41').Value = 'Patrik Udeh'
$ws.Range('K41').Value = 'Eric Wiberg'

$ws.Range('H42').Value = 0
$ws.Range('I42').Value = 1
$ws.Range('J42').Value = 'Dean Zoulamis'
$ws.Range('K42').Value = 'Nick Benson'

$ws.Range('H43').Value = 3
$ws.Range('I43').Value = 4
$ws.Range('J43').Value = 'Nick Benson'
$ws.Range('K43').Value = 'Dean Zoulamis'

$ws.Range('H44').Value = 5
$ws.Range('I44').Value = 2
$ws.Range('J44').Value = 'Eric Wiberg'
$ws.Range('K44').Value = 'Patrik Udeh'

$ws.Range('H45').Value = 2
$ws.Range('I45').Value = 1
$ws.Range('J45').Value = 'Patrik Udeh'
$ws.Range('K45').Value = 'Eric Wiberg'

$ws.Range('H46').Value = 1
$ws.Range('I46').Value = 0
$ws.Range('J46').Value = 'Dean Zoulamis'
$ws.Range('K46').Value = 'Patrik Udeh'

$ws.Range('H47').Value = 2
$ws.Range('I47').Value = 1
$ws.Range('J47').Value = 'Patrik Udeh'
$ws.Range('K47').Value = 'Dean Zoulamis'

$ws.Range('H48').Value = 0
$ws.Range('I48').Value = 7
$ws.Range('J48').Value = 'Eric Wiberg'
$ws.Range('K48').Value = 'Patrik Udeh'

$ws.Range('H49').Value = 3
$ws.Range('I49').Value = 3
$ws.Range('J49').Value = 'Patrik Udeh'
$ws.Range('K49').Value = 'Eric Wiberg'

$ws.Range('H50').Value = 4
$ws.Range('I50').Value = 1
$ws.Range('J50').Value = 'Nick Benson'
$ws.Range('K50').Value = 'Patrik Udeh'

$ws.Range('H51').Value = 0
$ws.Range('I51').Value = 1
$ws.Range('J51').Value = 'Patrik Udeh'
$ws.Range('K51').Value = 'Nick Benson'

$ws.Range('H52').Value = 1
$ws.Range('I52').Value = 2
$ws.Range('J52').Value = 'Eric Wiberg'
$ws.Range('K52').Value = 'Dean Zoulamis'

$ws.Range('H53').Value = 8
$ws.Range('I53').Value = 1
$ws.Range('J53').Value = 'Dean Zoulamis'
$ws.Range('K53').Value = 'Eric Wiberg'

$ws.Range('H54').Value = 0
$ws.Range('I54').Value = 1
$ws.Range('J54').Value = 'Nick Benson'
$ws.Range('K54').Value = 'Patrik Udeh'

$ws.Range('H55').Value = 1
$ws.Range('I55').Value = 2
$ws.Range('J55').Value = 'Patrik Udeh'
$ws.Range('K55').Value = 'Nick Benson'

$ws.Range('H56').Value = 0
$ws.Range('I56').Value = 2
$ws.Range('J56').Value = 'Dean Zoulamis'
$ws.Range('K56').Value = 'Eric Wiberg'

$ws.Range('H57').Value = 2
$ws.Range('I57').Value = 0
$ws.Range('J57').Value = 'Eric Wiberg'
$ws.Range('K57').Value = 'Dean Zoulamis'

$ws.Range('A58').Value = 'Quarter-Finals'
$ws.Range('B58').Value = 'October'
$ws.Range('C58').Value = 1
$ws.Range('D58').Value = 'Tuesday'
$ws.Range('E58').Value = 0.80902777777777779
$ws.Range('F58').Value = 'Thunder FC'
$ws.Range('G58').Value = 'APEX Charters Lone Pine Brewing'
$ws.Range('H58').Value = 7
$ws.Range('I58').Value = 2
$ws.Range('J58').Value = 'Patrik Udeh'
$ws.Range('K58').Value = 'Nick Benson'

$ws.Range('A59').Value = 'Quarter-Finals'
$ws.Range('B59').Value = 'October'
$ws.Range('C59').Value = 1
$ws.Range('D59').Value = 'Tuesday'
$ws.Range('E59').Value = 0.88194444444444453
$ws.Range('F59').Value = 'One Love FC'
$ws.Range('G59').Value = 'The Escape Room'
$ws.Range('H59').Value = 3
$ws.Range('I59').Value = 2
$ws.Range('J59').Value = 'Nick Benson'
$ws.Range('K59').Value = 'Patrik Udeh'

$ws.Range('A60').Value = 'Quarter-Finals'
$ws.Range('B60').Value = 'October'
$ws.Range('C60').Value = 3
$ws.Range('D60').Value = 'Thursday'
$ws.Range('E60').Value = 0.80902777777777779
$ws.Range('F60').Value = 'Carlos Auto Repair'
$ws.Range('G60').Value = 'Old Port FC'
$ws.Range('H60').Value = 1
$ws.Range('I60').Value = 0
$ws.Range('J60').Value = 'Patrik Udeh'
$ws.Range('K60').Value = 'Eric Wiberg'
$ws.Range('L60').Value = 'Ben (Pyta) Lomeri'

$ws.Range('A61').Value = 'Quarter-Finals'
$ws.Range('B61').Value = 'October'
$ws.Range('C61').Value = 3
$ws.Range('D61').Value = 'Thursday'
$ws.Range('E61').Value = 0.88194444444444453
$ws.Range('F61').Value = 'Farmers FC'
$ws.Range('G61').Value = 'Baxter Pines FC'
$ws.Range('H61').Value = 1
$ws.Range('I61').Value = 2
$ws.Range('J61').Value = 'Eric Wiberg'
$ws.Range('K61').Value = 'Patrik Udeh'
$ws.Range('L61').Value = 'Ben (Pyta) Lomeri'

$ws.Range('A62').Value = 'Semi-Finals'
$ws.Range('B62').Value = 'October'
$ws.Range('C62').Value = 8
$ws.Range('D62').Value = 'Tuesday'
$ws.Range('E62').Value = 0.80902777777777779
$ws.Range('F62').Value = 'Thunder FC'
$ws.Range('G62').Value = 'Carlos Auto Repair'
$ws.Range('H62').Value = 4
$ws.Range('I62').Value = 3
$ws.Range('J62').Value = 'Nick Benson'
$ws.Range('K62').Value = 'Eric Wiberg'
$ws.Range('L62').Value = 'Patrik Udeh'

$ws.Range('A63').Value = 'Semi-Finals'
$ws.Range('B63').Value = 'October'
$ws.Range('C63').Value = 8
$ws.Range('D63').Value = 'Tuesday'
$ws.Range('E63').Value = 0.88194444444444453
$ws.Range('F63').Value = 'Baxter Pines FC'
$ws.Range('G63').Value = 'One Love FC'
$ws.Range('H63').Value = 0
$ws.Range('I63').Value = 4
$ws.Range('J63').Value = 'Eric Wiberg'
$ws.Range('K63').Value = 'Nick Benson'
$ws.Range('L63').Value = 'Patrik Udeh'

$ws.Range('A64').Value = 'Final'
$ws.Range('B64').Value = 'October'
$ws.Range('C64').Value = 15
$ws.Range('D64').Value = 'Tuesday'
$ws.Range('E64').Value = 0.80902777777777779
$ws.Range('F64').Value = 'Thunder FC'
$ws.Range('G64').Value = 'One Love FC'
$ws.Range('H64').Value = 2
$ws.Range('I64').Value = 0
$ws.Range('J64').Value = 'Patrik Udeh'
$ws.Range('K64').Value = 'Dean Zoulamis'
$ws.Range('L64').Value = 'Eric Wiberg'

# --- Sheet view / selection updates ---
$ws4 = $wb.Worksheets.Item(4)  # 2023-Match

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range('K53').Select() | Out-Null

$ws4.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws4.Range('A58:A64').Select() | Out-Null
$ws.Activate()
